# Codon usage table: switch RNA-style codons (U) to DNA-style codons (T).
# Only column A (the "Codon" column) holds the letters that need the swap;
# column B (amino acid) and column C (frequency) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 65
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $codon = $cell.Value2
    if ($codon -ne $null) {
        $newCodon = $codon.Replace("U", "T")
        if ($newCodon -ne $codon) {
            $cell.Value = $newCodon
        }
    }
}

# Mirror the author's final UI state: column A ends up selected.
$null = $ws.Columns.Item(1).Select()
